$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 9 ("Diameter" row for the Directed Graphs section): 12 -> 14, 11 -> 12 ---

# South Africa column: "12" -> "14"
$cell = $t.Cell(9, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "14"

# Kenya column: "11" -> "12"
$cell = $t.Cell(9, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "12"

# --- Row 10 (new "Connected?" row) ---

# Label cell: bold "Connected?"
$cell = $t.Cell(10, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "Connected?"
$r.Font.Bold = $true
# Paragraph mark formatting also carries bold, matching the diff's pPr/rPr change.
$cell.Range.Paragraphs.Item(1).Range.Font.Bold = $true

# South Africa column: "No"
$cell = $t.Cell(10, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "No"

# Kenya column: "No"
$cell = $t.Cell(10, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "No"
